$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.010789
$ws.Range("H2").Value = 0.032367
$ws.Range("I2").Value = 0.01148982983039127
$ws.Range("J2").Value = 0.01148982983039127
$ws.Range("M2").Value = 0.003643333333333333
$ws.Range("N2").Value = 0.01093
$ws.Range("O2").Value = 0.002177035403614994
$ws.Range("P2").Value = 0.002177035403614994
$ws.Range("Q2").Value = 0.00003930792333333334
$ws.Range("R2").Value = 0.00035377131
$ws.Range("S2").Value = 0.00002501376632227346
$ws.Range("T2").Value = 0.00002501376632227346

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.010789
$ws.Range("H3").Value = 0.032367
$ws.Range("I3").Value = 0.01148982983039127
$ws.Range("J3").Value = 0.01148982983039127
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.669886333333333
$ws.Range("N3").Value = 5.009659
$ws.Range("O3").Value = 0.997822964596385
$ws.Range("P3").Value = 0.997822964596385
$ws.Range("Q3").Value = 0.01801640365033333
$ws.Range("R3").Value = 0.162147632853
$ws.Range("S3").Value = 0.01146481606406899
$ws.Range("T3").Value = 0.011464816064069

# Row 4
$ws.Range("I4").Value = 0.8731393855832401
$ws.Range("J4").Value = 0.8731393855832402
$ws.Range("M4").Value = 0.003643333333333333
$ws.Range("N4").Value = 0.01093
$ws.Range("O4").Value = 0.002177035403614994
$ws.Range("P4").Value = 0.002177035403614994
$ws.Range("Q4").Value = 0.002987102205555556
$ws.Range("R4").Value = 0.02688391985
$ws.Range("S4").Value = 0.001900855354705357
$ws.Range("T4").Value = 0.001900855354705357

# Row 5
$ws.Range("I5").Value = 0.8731393855832401
$ws.Range("J5").Value = 0.8731393855832402
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.669886333333333
$ws.Range("N5").Value = 5.009659
$ws.Range("O5").Value = 0.997822964596385
$ws.Range("P5").Value = 0.997822964596385
$ws.Range("Q5").Value = 1.369109190117222
$ws.Range("R5").Value = 12.321982711055
$ws.Range("S5").Value = 0.8712385302285347
$ws.Range("T5").Value = 0.8712385302285348

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.1083336666666667
$ws.Range("H6").Value = 0.325001
$ws.Range("I6").Value = 0.1153707845863686
$ws.Range("J6").Value = 0.1153707845863686
$ws.Range("M6").Value = 0.003643333333333333
$ws.Range("N6").Value = 0.01093
$ws.Range("O6").Value = 0.002177035403614994
$ws.Range("P6").Value = 0.002177035403614994
$ws.Range("Q6").Value = 0.0003946956588888889
$ws.Range("R6").Value = 0.00355226093
$ws.Range("S6").Value = 0.0002511662825873635
$ws.Range("T6").Value = 0.0002511662825873635

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.1083336666666667
$ws.Range("H7").Value = 0.325001
$ws.Range("I7").Value = 0.1153707845863686
$ws.Range("J7").Value = 0.1153707845863686
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.669886333333333
$ws.Range("N7").Value = 5.009659
$ws.Range("O7").Value = 0.997822964596385
$ws.Range("P7").Value = 0.997822964596385
$ws.Range("Q7").Value = 0.1809049094065555
$ws.Range("R7").Value = 1.628144184659
$ws.Range("S7").Value = 0.1151196183037812
$ws.Range("T7").Value = 0.1151196183037812
